# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets (zh-cn, de-de)
# now that a handback has happened, flips the Status from "Ready for
# handoff" to "Handed back: in sync with en-US" everywhere it appears, and
# widens a few columns so the newly-populated long file names are legible.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$githubBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/7ed0b58048952ab709cd959dea829061846689bb/e2e/"

$targetMdBe89 = "be89b7da-88b6-4d30-a0f4-ec864ee83e08.md"
$targetMdE87b = "e87b3ebd-18b1-49cf-b232-fe0371daea31.md"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (E, F) for both rows go from
# "Ready for handoff" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns to fit the longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = "be89b7da-88b6-4d30-a0f4-ec864ee83e08.4159531b2e61e89a33d76c47dedcdbe97d913bc5.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-12 18:59:17"

$wsZh.Range("J3").Value = "e87b3ebd-18b1-49cf-b232-fe0371daea31.778b0ea6a331d4b98101038488649d48da881242.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-12 18:59:17"

# Recreate the hyperlinks so the new "Latest Target File" links (I2, I3)
# sit alongside the existing Source File Name links (A2, A3), each
# pointing at the matching markdown file on GitHub.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $githubBase + $targetMdBe89, "", "", $targetMdBe89)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $githubBase + $targetMdBe89, "", "", $targetMdBe89)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $githubBase + $targetMdE87b, "", "", $targetMdE87b)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $githubBase + $targetMdE87b, "", "", $targetMdE87b)

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.16
$wsZh.Columns.Item(10).ColumnWidth = 39.16

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = "be89b7da-88b6-4d30-a0f4-ec864ee83e08.4159531b2e61e89a33d76c47dedcdbe97d913bc5.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-12 18:59:27"

$wsDe.Range("J3").Value = "e87b3ebd-18b1-49cf-b232-fe0371daea31.778b0ea6a331d4b98101038488649d48da881242.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-12 18:59:27"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $githubBase + $targetMdBe89, "", "", $targetMdBe89)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $githubBase + $targetMdBe89, "", "", $targetMdBe89)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $githubBase + $targetMdE87b, "", "", $targetMdE87b)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $githubBase + $targetMdE87b, "", "", $targetMdE87b)

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.16
$wsDe.Columns.Item(10).ColumnWidth = 39.16
